# Append: 2025-09-15 12:36 JST
# Refresh the 案件情報 (job listing) sheet with a new scrape snapshot:
#  - every existing row gets the new "取得日時" (fetched-at) timestamp
#  - 3 brand new listings are merged in at their appropriate positions
#  - hyperlinks in column F are rebuilt to match the rebuilt rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Remove all existing hyperlinks (and their relationships) before rewriting
# the row data; they will be recreated below for the new row layout.
$ws.Hyperlinks.Delete()

$timestamp = "2025-09-15 12:36:10"

$rows = @(
    @("【AI活用】データ分析Webサービス開発パートナー募集", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5393779", 368, "🔥AI,Ai ◆開発"),
    @("【急募】メモリデータ管理ツール開発のプロフェッショナル募集", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5393508", 158, "◆ツール,開発 ◇管理"),
    @("【簡易開発】会計・受発注システム付きITツールの依頼", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5393712", 143, "◆ツール,開発"),
    @("【急募】屋上貸切露天風呂の空き状況確認システム開発", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5389645", 125, "◆開発,システム開発"),
    @("FBA商品(在庫過多商品)をヤフオクで併売するツール開発依頼", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5393539", 123, "◆ツール,開発"),
    @("【急募】Gasを使用した公式LINEチャットbotの作成依頼", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5393641", 118, "★bot"),
    @("仮想通貨トレードの運用とコンサル【1名】のみ募集", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5393695", 55, "◆コンサル"),
    @("【急募】WordPressサーバー保守の専門家を探しています!(Xserver)", "システム開発", "5,000 円 ~ 10,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5393759", 25, "○WordPress"),
    @("【急募】トライアスロン大会運営支援システムの動作チェック、デバグ、品質確認業務委託費", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5393606", 33, $null)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $row[4])
    $ws.Cells.Item($r, 6).Style = "Hyperlink"
    $ws.Cells.Item($r, 7).Value = $row[5]
    if ($row[6] -ne $null) {
        $ws.Cells.Item($r, 8).Value = $row[6]
    }
    $r = $r + 1
}
